# Auto-generated Excel COM-interop script applying the Aegis_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 3781.375
$ws.Cells.Item(88, 9).Value = 1900.75
$ws.Cells.Item(88, 10).Value = 4408.25
$ws.Cells.Item(88, 11).Value = 1900.75
$ws.Cells.Item(88, 12).Value = 4408.25
$ws.Cells.Item(88, 13).Value = -1494.75
$ws.Cells.Item(88, 14).Value = -5220.25
$ws.Cells.Item(91, 8).Value = 3781.375
$ws.Cells.Item(91, 9).Value = 1900.75
$ws.Cells.Item(91, 10).Value = 4408.25
$ws.Cells.Item(91, 11).Value = 1900.75
$ws.Cells.Item(91, 12).Value = 4408.25
$ws.Cells.Item(91, 13).Value = -496.75
$ws.Cells.Item(91, 14).Value = -7216.25
$ws.Cells.Item(106, 8).Value = 2275.3
$ws.Cells.Item(106, 9).Value = 2361.4443
$ws.Cells.Item(106, 11).Value = 2361.4443
$ws.Cells.Item(106, 13).Value = -1730.4443
$ws.Cells.Item(113, 8).Value = 45400.26
$ws.Cells.Item(113, 9).Value = 92855
$ws.Cells.Item(113, 11).Value = 92855
$ws.Cells.Item(113, 13).Value = -89601
$ws.Cells.Item(135, 8).Value = 571.7778
$ws.Cells.Item(135, 9).Value = 586.5
$ws.Cells.Item(135, 10).Value = 520.25
$ws.Cells.Item(135, 11).Value = 5278.5
$ws.Cells.Item(135, 12).Value = 4682.25
$ws.Cells.Item(135, 13).Value = -2743.5
$ws.Cells.Item(135, 14).Value = -9752.25
$ws.Cells.Item(137, 8).Value = 1341.836
$ws.Cells.Item(137, 9).Value = 1064.8572
$ws.Cells.Item(137, 10).Value = 1576.8485
$ws.Cells.Item(137, 11).Value = 3194.5716
$ws.Cells.Item(137, 12).Value = 4730.5455
$ws.Cells.Item(137, 13).Value = -644.5715999999998
$ws.Cells.Item(137, 14).Value = -9830.5455
$ws.Cells.Item(139, 8).Value = 50780
$ws.Cells.Item(139, 10).Value = 50780
$ws.Cells.Item(139, 12).Value = 50780
$ws.Cells.Item(139, 14).Value = -61060

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1517.7091
$ws.Cells.Item(61, 9).Value = 689.63635
$ws.Cells.Item(61, 10).Value = 2069.7576
$ws.Cells.Item(61, 11).Value = 689.63635
$ws.Cells.Item(61, 12).Value = 2069.7576
$ws.Cells.Item(61, 13).Value = -477.63635
$ws.Cells.Item(61, 14).Value = -2493.7576
$ws.Cells.Item(74, 8).Value = 3057.1785
$ws.Cells.Item(74, 9).Value = 2490.5454
$ws.Cells.Item(74, 10).Value = 3423.8235
$ws.Cells.Item(74, 11).Value = 2490.5454
$ws.Cells.Item(74, 12).Value = 3423.8235
$ws.Cells.Item(74, 13).Value = -1616.5454
$ws.Cells.Item(74, 14).Value = -5171.8235
$ws.Cells.Item(77, 8).Value = 3057.1785
$ws.Cells.Item(77, 9).Value = 2490.5454
$ws.Cells.Item(77, 10).Value = 3423.8235
$ws.Cells.Item(77, 11).Value = 12452.727
$ws.Cells.Item(77, 12).Value = 17119.1175
$ws.Cells.Item(77, 13).Value = -8084.726999999999
$ws.Cells.Item(77, 14).Value = -25855.1175
$ws.Cells.Item(110, 8).Value = 45501876
$ws.Cells.Item(110, 9).Value = 77001224
$ws.Cells.Item(110, 10).Value = 2812.4443
$ws.Cells.Item(110, 11).Value = 77001224
$ws.Cells.Item(110, 12).Value = 2812.4443
$ws.Cells.Item(110, 13).Value = -76999179
$ws.Cells.Item(110, 14).Value = -6902.4443
$ws.Cells.Item(132, 8).Value = 2960.5098
$ws.Cells.Item(132, 9).Value = 3775.1
$ws.Cells.Item(132, 10).Value = 1796.8096
$ws.Cells.Item(132, 11).Value = 11325.3
$ws.Cells.Item(132, 12).Value = 5390.4288
$ws.Cells.Item(132, 13).Value = -8795.299999999999
$ws.Cells.Item(132, 14).Value = -10450.4288
$ws.Cells.Item(136, 8).Value = 1517.7091
$ws.Cells.Item(136, 9).Value = 689.63635
$ws.Cells.Item(136, 10).Value = 2069.7576
$ws.Cells.Item(136, 11).Value = 2068.90905
$ws.Cells.Item(136, 12).Value = 6209.2728
$ws.Cells.Item(136, 13).Value = 481.0909499999998
$ws.Cells.Item(136, 14).Value = -11309.2728
$ws.Cells.Item(140, 8).Value = 65424.5
$ws.Cells.Item(140, 10).Value = 65424.5
$ws.Cells.Item(140, 12).Value = 65424.5
$ws.Cells.Item(140, 14).Value = -75784.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(34, 8).Value = 49800
$ws.Cells.Item(34, 10).Value = 49800
$ws.Cells.Item(34, 12).Value = 49800
$ws.Cells.Item(34, 14).Value = -50028
$ws.Cells.Item(35, 8).Value = 19322.8
$ws.Cells.Item(35, 10).Value = 19322.8
$ws.Cells.Item(35, 12).Value = 19322.8
$ws.Cells.Item(35, 14).Value = -19942.8
$ws.Cells.Item(94, 8).Value = 143212.14
$ws.Cells.Item(94, 9).Value = 143212.14
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 143212.14
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -142761.14
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 3135.55
$ws.Cells.Item(134, 9).Value = 3206.8125
$ws.Cells.Item(134, 11).Value = 9620.4375
$ws.Cells.Item(134, 13).Value = -7085.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 46561.332
$ws.Cells.Item(20, 10).Value = 46561.332
$ws.Cells.Item(20, 12).Value = 46561.332
$ws.Cells.Item(20, 14).Value = -47033.332
$ws.Cells.Item(30, 8).Value = 46561.332
$ws.Cells.Item(30, 10).Value = 46561.332
$ws.Cells.Item(30, 12).Value = 46561.332
$ws.Cells.Item(30, 14).Value = -46743.332
$ws.Cells.Item(31, 8).Value = 14471.684
$ws.Cells.Item(31, 9).Value = 29663.771
$ws.Cells.Item(31, 10).Value = 2387.068
$ws.Cells.Item(31, 11).Value = 29663.771
$ws.Cells.Item(31, 12).Value = 2387.068
$ws.Cells.Item(31, 13).Value = -29368.771
$ws.Cells.Item(31, 14).Value = -2977.068
$ws.Cells.Item(34, 8).Value = 14471.684
$ws.Cells.Item(34, 9).Value = 29663.771
$ws.Cells.Item(34, 10).Value = 2387.068
$ws.Cells.Item(34, 11).Value = 29663.771
$ws.Cells.Item(34, 12).Value = 2387.068
$ws.Cells.Item(34, 13).Value = -29461.771
$ws.Cells.Item(34, 14).Value = -2791.068
$ws.Cells.Item(58, 8).Value = 12265.667
$ws.Cells.Item(58, 9).Value = 1598
$ws.Cells.Item(58, 10).Value = 52802.8
$ws.Cells.Item(58, 11).Value = 1598
$ws.Cells.Item(58, 12).Value = 52802.8
$ws.Cells.Item(58, 13).Value = -1395
$ws.Cells.Item(58, 14).Value = -53208.8
$ws.Cells.Item(107, 8).Value = 654.3333
$ws.Cells.Item(107, 9).Value = 623.2778
$ws.Cells.Item(107, 11).Value = 623.2778
$ws.Cells.Item(107, 13).Value = 1296.7222
$ws.Cells.Item(128, 8).Value = 46561.332
$ws.Cells.Item(128, 10).Value = 46561.332
$ws.Cells.Item(128, 12).Value = 46561.332
$ws.Cells.Item(128, 14).Value = -56521.332
$ws.Cells.Item(136, 8).Value = 12265.667
$ws.Cells.Item(136, 9).Value = 1598
$ws.Cells.Item(136, 10).Value = 52802.8
$ws.Cells.Item(136, 11).Value = 4794
$ws.Cells.Item(136, 12).Value = 158408.4
$ws.Cells.Item(136, 13).Value = -2244
$ws.Cells.Item(136, 14).Value = -163508.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 53.77778
$ws.Cells.Item(2, 9).Value = 19.75
$ws.Cells.Item(2, 10).Value = 81
$ws.Cells.Item(2, 11).Value = 118.5
$ws.Cells.Item(2, 12).Value = 486
$ws.Cells.Item(2, 13).Value = -5.5
$ws.Cells.Item(2, 14).Value = -712
$ws.Cells.Item(17, 8).Value = 2142.8572
$ws.Cells.Item(17, 9).Value = 2000
$ws.Cells.Item(17, 10).Value = 3000
$ws.Cells.Item(17, 11).Value = 6000
$ws.Cells.Item(17, 12).Value = 9000
$ws.Cells.Item(17, 13).Value = -5831
$ws.Cells.Item(17, 14).Value = -9338
$ws.Cells.Item(34, 8).Value = 1659
$ws.Cells.Item(34, 10).Value = 2237.375
$ws.Cells.Item(34, 12).Value = 6712.125
$ws.Cells.Item(34, 14).Value = -6880.125
$ws.Cells.Item(39, 8).Value = 2980
$ws.Cells.Item(39, 10).Value = 2980
$ws.Cells.Item(39, 12).Value = 8940
$ws.Cells.Item(39, 14).Value = -9528
$ws.Cells.Item(55, 8).Value = 8279.5
$ws.Cells.Item(55, 10).Value = 8279.5
$ws.Cells.Item(55, 12).Value = 24838.5
$ws.Cells.Item(55, 14).Value = -25192.5
$ws.Cells.Item(68, 8).Value = 2015.3151
$ws.Cells.Item(68, 9).Value = 1332.3549
$ws.Cells.Item(68, 10).Value = 2519.4048
$ws.Cells.Item(68, 11).Value = 3997.0647
$ws.Cells.Item(68, 12).Value = 7558.214399999999
$ws.Cells.Item(68, 13).Value = -3186.0647
$ws.Cells.Item(68, 14).Value = -9180.214399999999
$ws.Cells.Item(71, 8).Value = 2015.3151
$ws.Cells.Item(71, 9).Value = 1332.3549
$ws.Cells.Item(71, 10).Value = 2519.4048
$ws.Cells.Item(71, 11).Value = 11991.1941
$ws.Cells.Item(71, 12).Value = 22674.6432
$ws.Cells.Item(71, 13).Value = -7935.194100000001
$ws.Cells.Item(71, 14).Value = -30786.6432
$ws.Cells.Item(131, 8).Value = 1334786.4
$ws.Cells.Item(131, 9).Value = 1177
$ws.Cells.Item(131, 10).Value = 1551885.5
$ws.Cells.Item(131, 11).Value = 3531
$ws.Cells.Item(131, 12).Value = 4655656.5
$ws.Cells.Item(131, 13).Value = 1509
$ws.Cells.Item(131, 14).Value = -4665736.5
$ws.Cells.Item(137, 8).Value = 2573.8696
$ws.Cells.Item(137, 10).Value = 2854.2144
$ws.Cells.Item(137, 12).Value = 8562.643199999999
$ws.Cells.Item(137, 14).Value = -18762.6432
$ws.Cells.Item(140, 8).Value = 1611.16
$ws.Cells.Item(140, 9).Value = 1155
$ws.Cells.Item(140, 10).Value = 2422.111
$ws.Cells.Item(140, 11).Value = 3465
$ws.Cells.Item(140, 12).Value = 7266.333
$ws.Cells.Item(140, 13).Value = 1715
$ws.Cells.Item(140, 14).Value = -17626.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2231.3845
$ws.Cells.Item(126, 9).Value = 2372.5715
$ws.Cells.Item(126, 11).Value = 7117.7145
$ws.Cells.Item(126, 13).Value = -4647.7145
$ws.Cells.Item(132, 8).Value = 2044.4166
$ws.Cells.Item(132, 9).Value = 1469.4
$ws.Cells.Item(132, 10).Value = 3002.7778
$ws.Cells.Item(132, 11).Value = 4408.200000000001
$ws.Cells.Item(132, 12).Value = 9008.3334
$ws.Cells.Item(132, 13).Value = -1878.200000000001
$ws.Cells.Item(132, 14).Value = -14068.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 5227.273
$ws.Cells.Item(2, 9).Value = 1500
$ws.Cells.Item(2, 10).Value = 8333.333000000001
$ws.Cells.Item(2, 11).Value = 1500
$ws.Cells.Item(2, 12).Value = 8333.333000000001
$ws.Cells.Item(2, 13).Value = -1388
$ws.Cells.Item(2, 14).Value = -8557.333000000001
$ws.Cells.Item(122, 8).Value = 5000
$ws.Cells.Item(122, 9).Value = 5000
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 15000
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -12550
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 4889.85
$ws.Cells.Item(132, 9).Value = 6118.364
$ws.Cells.Item(132, 10).Value = 3388.3333
$ws.Cells.Item(132, 11).Value = 18355.092
$ws.Cells.Item(132, 12).Value = 10164.9999
$ws.Cells.Item(132, 13).Value = -15825.092
$ws.Cells.Item(132, 14).Value = -15224.9999
$ws.Cells.Item(136, 8).Value = 1941.4615
$ws.Cells.Item(136, 9).Value = 1587.2727
$ws.Cells.Item(136, 10).Value = 3889.5
$ws.Cells.Item(136, 11).Value = 4761.8181
$ws.Cells.Item(136, 12).Value = 11668.5
$ws.Cells.Item(136, 13).Value = -2211.8181
$ws.Cells.Item(136, 14).Value = -16768.5
$ws.Cells.Item(140, 8).Value = 48583.855
$ws.Cells.Item(140, 10).Value = 48583.855
$ws.Cells.Item(140, 12).Value = 48583.855
$ws.Cells.Item(140, 14).Value = -58943.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 2505000
$ws.Cells.Item(5, 10).Value = 5000000
$ws.Cells.Item(5, 12).Value = 5000000
$ws.Cells.Item(5, 14).Value = -5000224
$ws.Cells.Item(16, 8).Value = 40233.332
$ws.Cells.Item(16, 9).Value = 38700
$ws.Cells.Item(16, 10).Value = 41000
$ws.Cells.Item(16, 11).Value = 38700
$ws.Cells.Item(16, 12).Value = 41000
$ws.Cells.Item(16, 13).Value = -38408
$ws.Cells.Item(16, 14).Value = -41584
$ws.Cells.Item(96, 8).Value = 250002100
$ws.Cells.Item(96, 9).Value = 333335330
$ws.Cells.Item(96, 10).Value = 2404
$ws.Cells.Item(96, 11).Value = 333335330
$ws.Cells.Item(96, 12).Value = 2404
$ws.Cells.Item(96, 13).Value = -333333957
$ws.Cells.Item(96, 14).Value = -5150
$ws.Cells.Item(132, 8).Value = 2905.3057
$ws.Cells.Item(132, 9).Value = 3476.9546
$ws.Cells.Item(132, 11).Value = 10430.8638
$ws.Cells.Item(132, 13).Value = -7900.863799999999

Write-Host "Applied all changes"